# "Generate Report for handoff"
# Updates the localization-status report: marks files as ready for handoff,
# records the generated handoff (.xlf) package links + timestamps, and
# flips the dependency handling from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/830b991ba7253bfbef1f6cf41e8a7635d8987da0/e2e/"

function Update-LangSheet($SheetName, $XlfName, $HandoffDatetime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # B2: Status -> "Ready for handoff"
    $ws.Range("B2").Value = "Ready for handoff"

    # C2: Latest Handoff File -> new .xlf package, as a hyperlink
    $ws.Range("C2").Value = $XlfName
    $ws.Hyperlinks.Add($ws.Range("C2"), ($baseUrl + $XlfName), "", "", $XlfName)
    $ws.Range("C2").Font.Underline = $true
    $ws.Range("C2").Font.Color = [System.Convert]::ToInt64("ED9564", 16)

    # D2: Latest Handoff Datetime -> generation time of the handoff package
    $ws.Range("D2").Value = $HandoffDatetime
    $ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    # H2: Handoff Reason -> "Include" (was "Ignored")
    $ws.Range("H2").Value = "Include"
}

Update-LangSheet "zh-cn" "f939a530-8a03-4901-891d-bcc658750a13.5feeb2162ecc4446a2c0e819b46d3eea54254d29.zh-cn.xlf" "2016-01-13 12:58:30"
Update-LangSheet "de-de" "f939a530-8a03-4901-891d-bcc658750a13.5feeb2162ecc4446a2c0e819b46d3eea54254d29.de-de.xlf" "2016-01-13 12:58:55"

# The Overview sheet mirrors each language sheet's Status cell (shares the
# same "Handoff transform failed" -> "Ready for handoff" text), so it needs
# the same update.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
